$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header (default / primary) ---
$hdr = $sec.Headers.Item(1)
$hdr.Range.Text = "27/05/2024"
$hdr.Range.Font.Name = "Arial"
$hdr.Range.Font.Bold = $true
$hdr.Range.ParagraphFormat.Alignment = 2

# --- Footer (default / primary) ---
$ftr = $sec.Footers.Item(1)
$ftr.Range.Text = "4IV9 Quirino González Johann David"
$ftr.Range.Font.Name = "Arial"
$ftr.Range.Font.Bold = $true
$ftr.Range.ParagraphFormat.Alignment = 2

Write-Output "done"
